$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REST-API-LIST")
$ws.Range("A30").Style = "Bad"
$rst = $ws.Range("A30").Style
Write-Output "Type: $($rst.GetType())"
Write-Output "Name: $($rst.Name)"
$rst.Font.Name = "SimSun"
Write-Output "after set"
